$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Domain")

# Update C6 text (Patient Transfer description): "shifting" -> "transferring"
$ws.Range("C6").Value = "The use case deals with transferring the patient to an outpatient clinic based on the criticality of condition, entitlement of coverage."

# Update the view: scroll so column C is the left-most visible column,
# and select cell C11
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("C11").Select()
